$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 87500
$ws.Range("I13").Value = 75000
$ws.Range("J13").Value = 100000
$ws.Range("K13").Value = 75000
$ws.Range("L13").Value = 100000
$ws.Range("M13").Value = -74831
$ws.Range("N13").Value = -100338
$ws.Range("H17").Value = 581092.25
$ws.Range("J17").Value = 977000.5
$ws.Range("L17").Value = 2931001.5
$ws.Range("N17").Value = -2931337.5
$ws.Range("H34").Value = 2163.4
$ws.Range("I34").Value = 2163.4
$ws.Range("K34").Value = 2163.4
$ws.Range("M34").Value = -1960.4
$ws.Range("H36").Value = 2163.4
$ws.Range("I36").Value = 2163.4
$ws.Range("K36").Value = 2163.4
$ws.Range("M36").Value = -1448.4
$ws.Range("H40").Value = 2570.7273
$ws.Range("I40").Value = 2357.8
$ws.Range("J40").Value = 2748.1667
$ws.Range("K40").Value = 2357.8
$ws.Range("L40").Value = 2748.1667
$ws.Range("M40").Value = -2182.8
$ws.Range("N40").Value = -3098.1667
$ws.Range("H43").Value = 18866.75
$ws.Range("I43").Value = 26660
$ws.Range("J43").Value = 5878
$ws.Range("K43").Value = 26660
$ws.Range("L43").Value = 5878
$ws.Range("M43").Value = -26591
$ws.Range("N43").Value = -6016
$ws.Range("H51").Value = 2670.818
$ws.Range("I51").Value = 2313.1875
$ws.Range("J51").Value = 3624.5
$ws.Range("K51").Value = 2313.1875
$ws.Range("L51").Value = 3624.5
$ws.Range("M51").Value = -1829.1875
$ws.Range("N51").Value = -4592.5
$ws.Range("H70").Value = 55872.21
$ws.Range("J70").Value = 65966.875
$ws.Range("L70").Value = 197900.625
$ws.Range("N70").Value = -198440.625
$ws.Range("H73").Value = 55872.21
$ws.Range("J73").Value = 65966.875
$ws.Range("L73").Value = 197900.625
$ws.Range("N73").Value = -199772.625
$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("N88").Value = -2812
$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("N91").Value = -4808
$ws.Range("H107").Value = 362.42105
$ws.Range("I107").Value = 404.26666
$ws.Range("J107").Value = 205.5
$ws.Range("K107").Value = 404.26666
$ws.Range("L107").Value = 205.5
$ws.Range("M107").Value = 1515.73334
$ws.Range("N107").Value = -4045.5
$ws.Range("H130").Value = 154320
$ws.Range("J130").Value = 154320
$ws.Range("L130").Value = 154320
$ws.Range("N130").Value = -164360
$ws.Range("H131").Value = 40524.75
$ws.Range("I131").Value = 4995
$ws.Range("J131").Value = 52368
$ws.Range("K131").Value = 14985
$ws.Range("L131").Value = 157104
$ws.Range("M131").Value = -9945
$ws.Range("N131").Value = -167184
$ws.Range("H132").Value = 2468.7273
$ws.Range("I132").Value = 1962.9445
$ws.Range("J132").Value = 4744.75
$ws.Range("K132").Value = 5888.833500000001
$ws.Range("L132").Value = 14234.25
$ws.Range("M132").Value = -3358.833500000001
$ws.Range("N132").Value = -19294.25
$ws.Range("H138").Value = 11114829
$ws.Range("J138").Value = 15156100
$ws.Range("L138").Value = 45468300
$ws.Range("N138").Value = -45478580

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1947.8889
$ws.Range("I2").Value = 1947.8889
$ws.Range("K2").Value = 1947.8889
$ws.Range("M2").Value = -1834.8889
$ws.Range("H45").Value = 14258.556
$ws.Range("I45").Value = 14258.556
$ws.Range("K45").Value = 14258.556
$ws.Range("M45").Value = -13881.556
$ws.Range("H61").Value = 4333.0713
$ws.Range("I61").Value = 3113.1365
$ws.Range("K61").Value = 3113.1365
$ws.Range("M61").Value = -2901.1365
$ws.Range("H63").Value = 4614.4287
$ws.Range("I63").Value = 5261
$ws.Range("K63").Value = 5261
$ws.Range("M63").Value = -4575
$ws.Range("H64").Value = 39799.8
$ws.Range("I64").Value = 19999
$ws.Range("J64").Value = 44750
$ws.Range("K64").Value = 19999
$ws.Range("L64").Value = 44750
$ws.Range("M64").Value = -19751
$ws.Range("N64").Value = -45246
$ws.Range("H66").Value = 4614.4287
$ws.Range("I66").Value = 5261
$ws.Range("K66").Value = 26305
$ws.Range("M66").Value = -22873
$ws.Range("H67").Value = 39799.8
$ws.Range("I67").Value = 19999
$ws.Range("J67").Value = 44750
$ws.Range("K67").Value = 19999
$ws.Range("L67").Value = 44750
$ws.Range("M67").Value = -19141
$ws.Range("N67").Value = -46466
$ws.Range("H74").Value = 87098.30499999999
$ws.Range("I74").Value = 102817.8
$ws.Range("J74").Value = 34700
$ws.Range("K74").Value = 102817.8
$ws.Range("L74").Value = 34700
$ws.Range("M74").Value = -101943.8
$ws.Range("N74").Value = -36448
$ws.Range("H77").Value = 87098.30499999999
$ws.Range("I77").Value = 102817.8
$ws.Range("J77").Value = 34700
$ws.Range("K77").Value = 514089
$ws.Range("L77").Value = 173500
$ws.Range("M77").Value = -509721
$ws.Range("N77").Value = -182236
$ws.Range("H88").Value = 2621.5334
$ws.Range("I88").Value = 1777.5
$ws.Range("J88").Value = 2751.3845
$ws.Range("K88").Value = 1777.5
$ws.Range("L88").Value = 2751.3845
$ws.Range("M88").Value = -1371.5
$ws.Range("N88").Value = -3563.3845
$ws.Range("H91").Value = 2621.5334
$ws.Range("I91").Value = 1777.5
$ws.Range("J91").Value = 2751.3845
$ws.Range("K91").Value = 1777.5
$ws.Range("L91").Value = 2751.3845
$ws.Range("M91").Value = -373.5
$ws.Range("N91").Value = -5559.3845
$ws.Range("H102").Value = 2006.3914
$ws.Range("I102").Value = 1911.8948
$ws.Range("J102").Value = 2455.25
$ws.Range("K102").Value = 1911.8948
$ws.Range("L102").Value = 2455.25
$ws.Range("M102").Value = -289.8948
$ws.Range("N102").Value = -5699.25
$ws.Range("H116").Value = 1947.8889
$ws.Range("I116").Value = 1947.8889
$ws.Range("K116").Value = 1947.8889
$ws.Range("M116").Value = 346.1111000000001
$ws.Range("H122").Value = 2835.077
$ws.Range("I122").Value = 2762.889
$ws.Range("J122").Value = 2997.5
$ws.Range("K122").Value = 8288.667000000001
$ws.Range("L122").Value = 8992.5
$ws.Range("M122").Value = -5838.667000000001
$ws.Range("N122").Value = -13892.5
$ws.Range("H131").Value = 161109.5
$ws.Range("J131").Value = 161109.5
$ws.Range("L131").Value = 161109.5
$ws.Range("N131").Value = -171189.5
$ws.Range("H132").Value = 6921.7617
$ws.Range("I132").Value = 7085.5
$ws.Range("J132").Value = 6397.8
$ws.Range("K132").Value = 21256.5
$ws.Range("L132").Value = 19193.4
$ws.Range("M132").Value = -18726.5
$ws.Range("N132").Value = -24253.4
$ws.Range("H136").Value = 4333.0713
$ws.Range("I136").Value = 3113.1365
$ws.Range("K136").Value = 9339.4095
$ws.Range("M136").Value = -6789.4095

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1947.8889
$ws.Range("I3").Value = 1947.8889
$ws.Range("K3").Value = 1947.8889
$ws.Range("M3").Value = -1833.8889
$ws.Range("H99").Value = 5390.909
$ws.Range("I99").Value = 2144.4443
$ws.Range("K99").Value = 2144.4443
$ws.Range("M99").Value = -646.4443000000001
$ws.Range("H132").Value = 77600
$ws.Range("J132").Value = 77600
$ws.Range("L132").Value = 77600
$ws.Range("N132").Value = -87720
$ws.Range("H134").Value = 1914.807
$ws.Range("I134").Value = 1910.5962
$ws.Range("K134").Value = 5731.7886
$ws.Range("M134").Value = -3196.7886

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 926.93335
$ws.Range("I7").Value = 1157.1818
$ws.Range("K7").Value = 1157.1818
$ws.Range("M7").Value = -1044.1818
$ws.Range("H31").Value = 128426.625
$ws.Range("I31").Value = 253125
$ws.Range("J31").Value = 3728.25
$ws.Range("K31").Value = 253125
$ws.Range("L31").Value = 3728.25
$ws.Range("M31").Value = -252830
$ws.Range("N31").Value = -4318.25
$ws.Range("H34").Value = 128426.625
$ws.Range("I34").Value = 253125
$ws.Range("J34").Value = 3728.25
$ws.Range("K34").Value = 253125
$ws.Range("L34").Value = 3728.25
$ws.Range("M34").Value = -252923
$ws.Range("N34").Value = -4132.25
$ws.Range("H58").Value = 2777.3333
$ws.Range("I58").Value = 2714.1428
$ws.Range("K58").Value = 2714.1428
$ws.Range("M58").Value = -2511.1428
$ws.Range("H105").Value = 2365.625
$ws.Range("I105").Value = 1308.5
$ws.Range("K105").Value = 1308.5
$ws.Range("M105").Value = 438.5
$ws.Range("H122").Value = 2245.5
$ws.Range("I122").Value = 2025
$ws.Range("J122").Value = 2466
$ws.Range("K122").Value = 6075
$ws.Range("L122").Value = 7398
$ws.Range("M122").Value = -3625
$ws.Range("N122").Value = -12298
$ws.Range("H134").Value = 32570.857
$ws.Range("I134").Value = 10999.5
$ws.Range("J134").Value = 61332.668
$ws.Range("K134").Value = 32998.5
$ws.Range("L134").Value = 183998.004
$ws.Range("M134").Value = -30463.5
$ws.Range("N134").Value = -189068.004
$ws.Range("H136").Value = 2777.3333
$ws.Range("I136").Value = 2714.1428
$ws.Range("K136").Value = 8142.428400000001
$ws.Range("M136").Value = -5592.428400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34143.133
$ws.Range("I2").Value = 36574.93
$ws.Range("K2").Value = 219449.58
$ws.Range("M2").Value = -219336.58
$ws.Range("H15").Value = 23.6
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H18").Value = 931.5
$ws.Range("J18").Value = 1033
$ws.Range("L18").Value = 3099
$ws.Range("N18").Value = -3437
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H21").Value = 282.8889
$ws.Range("I21").Value = 215.85715
$ws.Range("J21").Value = 517.5
$ws.Range("K21").Value = 647.5714499999999
$ws.Range("L21").Value = 1552.5
$ws.Range("M21").Value = -474.5714499999999
$ws.Range("N21").Value = -1898.5
$ws.Range("H23").Value = 346.66666
$ws.Range("I23").Value = 283.8
$ws.Range("J23").Value = 378.1
$ws.Range("K23").Value = 851.4000000000001
$ws.Range("L23").Value = 1134.3
$ws.Range("M23").Value = -616.4000000000001
$ws.Range("N23").Value = -1604.3
$ws.Range("H38").Value = 67.73333
$ws.Range("I38").Value = 56.555557
$ws.Range("J38").Value = 84.5
$ws.Range("K38").Value = 169.666671
$ws.Range("L38").Value = 253.5
$ws.Range("M38").Value = 177.333329
$ws.Range("N38").Value = -947.5
$ws.Range("H92").Value = 539.7778
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 539.7778
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1619.3334
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -4115.3334
$ws.Range("H113").Value = 811.0769
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 813.0909
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 2439.2727
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -6779.2727

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18000000
$ws.Range("I11").Value = 18000000
$ws.Range("K11").Value = 18000000
$ws.Range("M11").Value = -17999861
$ws.Range("H18").Value = 55605556
$ws.Range("I18").Value = 111111110
$ws.Range("K18").Value = 111111110
$ws.Range("M18").Value = -111110817
$ws.Range("H21").Value = 11990
$ws.Range("I21").Value = 11990
$ws.Range("K21").Value = 11990
$ws.Range("M21").Value = -11817
$ws.Range("H30").Value = 11990
$ws.Range("I30").Value = 11990
$ws.Range("K30").Value = 11990
$ws.Range("M30").Value = -11885
$ws.Range("H70").Value = 19999
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 19999
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 19999
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -20539
$ws.Range("H73").Value = 19999
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 19999
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 19999
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -21871
$ws.Range("H80").Value = 2726.2727
$ws.Range("I80").Value = 2582.5
$ws.Range("J80").Value = 2898.8
$ws.Range("K80").Value = 2582.5
$ws.Range("L80").Value = 2898.8
$ws.Range("M80").Value = -1584.5
$ws.Range("N80").Value = -4894.8
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 30000
$ws.Range("K82").Value = 30000
$ws.Range("M82").Value = -29617
$ws.Range("H83").Value = 2726.2727
$ws.Range("I83").Value = 2582.5
$ws.Range("J83").Value = 2898.8
$ws.Range("K83").Value = 12912.5
$ws.Range("L83").Value = 14494
$ws.Range("M83").Value = -7920.5
$ws.Range("N83").Value = -24478
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 30000
$ws.Range("K85").Value = 30000
$ws.Range("M85").Value = -28674
$ws.Range("H97").Value = 1183.3334
$ws.Range("I97").Value = 1284.4445
$ws.Range("J97").Value = 1031.6666
$ws.Range("K97").Value = 1284.4445
$ws.Range("L97").Value = 1031.6666
$ws.Range("M97").Value = -788.4445000000001
$ws.Range("N97").Value = -2023.6666
$ws.Range("H113").Value = 40460.46
$ws.Range("I113").Value = 47498.727
$ws.Range("K113").Value = 47498.727
$ws.Range("M113").Value = -45328.727
$ws.Range("H119").Value = 29499.166
$ws.Range("J119").Value = 29499.166
$ws.Range("L119").Value = 29499.166
$ws.Range("N119").Value = -39175.166
$ws.Range("H122").Value = 3260.1875
$ws.Range("I122").Value = 3210.818
$ws.Range("K122").Value = 9632.454000000002
$ws.Range("M122").Value = -7182.454000000002
$ws.Range("H126").Value = 18446.6
$ws.Range("I126").Value = 22178.25
$ws.Range("J126").Value = 3520
$ws.Range("K126").Value = 66534.75
$ws.Range("L126").Value = 10560
$ws.Range("M126").Value = -64064.75
$ws.Range("N126").Value = -15500
$ws.Range("H132").Value = 3665.5334
$ws.Range("I132").Value = 2471.6365
$ws.Range("K132").Value = 7414.9095
$ws.Range("M132").Value = -4884.9095
$ws.Range("H136").Value = 32147.143
$ws.Range("J136").Value = 32147.143
$ws.Range("L136").Value = 96441.429
$ws.Range("N136").Value = -101541.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5717
$ws.Range("I40").Value = 5022.769
$ws.Range("K40").Value = 5022.769
$ws.Range("M40").Value = -4886.769
$ws.Range("H61").Value = 45607.5
$ws.Range("I61").Value = 49253.637
$ws.Range("K61").Value = 49253.637
$ws.Range("M61").Value = -49051.637
$ws.Range("H82").Value = 2591.5715
$ws.Range("I82").Value = 3012.625
$ws.Range("J82").Value = 2030.1666
$ws.Range("K82").Value = 3012.625
$ws.Range("L82").Value = 2030.1666
$ws.Range("M82").Value = -2651.625
$ws.Range("N82").Value = -2752.1666
$ws.Range("H85").Value = 2591.5715
$ws.Range("I85").Value = 3012.625
$ws.Range("J85").Value = 2030.1666
$ws.Range("K85").Value = 3012.625
$ws.Range("L85").Value = 2030.1666
$ws.Range("M85").Value = -1764.625
$ws.Range("N85").Value = -4526.1666
$ws.Range("H113").Value = 45607.5
$ws.Range("I113").Value = 49253.637
$ws.Range("K113").Value = 49253.637
$ws.Range("M113").Value = -47083.637
$ws.Range("H122").Value = 389365.8
$ws.Range("I122").Value = 530343.2
$ws.Range("J122").Value = 6713
$ws.Range("K122").Value = 1591029.6
$ws.Range("L122").Value = 20139
$ws.Range("M122").Value = -1588579.6
$ws.Range("N122").Value = -25039
$ws.Range("H136").Value = 4371.9565
$ws.Range("I136").Value = 4003.2354
$ws.Range("J136").Value = 5416.6665
$ws.Range("K136").Value = 12009.7062
$ws.Range("L136").Value = 16249.9995
$ws.Range("M136").Value = -9459.706200000001
$ws.Range("N136").Value = -21349.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 22200
$ws.Range("I32").Value = 20000
$ws.Range("K32").Value = 20000
$ws.Range("M32").Value = -19683
$ws.Range("H41").Value = 79251.336
$ws.Range("J41").Value = 79251.336
$ws.Range("L41").Value = 79251.336
$ws.Range("N41").Value = -80031.336
$ws.Range("H46").Value = 74285.8
$ws.Range("J46").Value = 74285.8
$ws.Range("L46").Value = 74285.8
$ws.Range("N46").Value = -74747.8
$ws.Range("H81").Value = 9203.049999999999
$ws.Range("J81").Value = 4808.8823
$ws.Range("L81").Value = 9617.7646
$ws.Range("N81").Value = -11739.7646
$ws.Range("H84").Value = 9203.049999999999
$ws.Range("J84").Value = 4808.8823
$ws.Range("L84").Value = 48088.823
$ws.Range("N84").Value = -58696.823
$ws.Range("H100").Value = 709.5625
$ws.Range("I100").Value = 724.6667
$ws.Range("J100").Value = 483
$ws.Range("K100").Value = 1449.3334
$ws.Range("L100").Value = 966
$ws.Range("M100").Value = -908.3334
$ws.Range("N100").Value = -2048
$ws.Range("H122").Value = 2452.4443
$ws.Range("I122").Value = 2082.9333
$ws.Range("J122").Value = 4300
$ws.Range("K122").Value = 6248.7999
$ws.Range("L122").Value = 12900
$ws.Range("M122").Value = -3798.7999
$ws.Range("N122").Value = -17800
$ws.Range("H132").Value = 1814
$ws.Range("I132").Value = 1662.0476
$ws.Range("J132").Value = 5005
$ws.Range("K132").Value = 4986.142800000001
$ws.Range("L132").Value = 15015
$ws.Range("M132").Value = -2456.142800000001
$ws.Range("N132").Value = -20075
$ws.Range("H134").Value = 74285.8
$ws.Range("J134").Value = 74285.8
$ws.Range("L134").Value = 222857.4
$ws.Range("N134").Value = -227927.4
$ws.Range("H136").Value = 2819.923
$ws.Range("I136").Value = 2696.0454
$ws.Range("K136").Value = 8088.1362
$ws.Range("M136").Value = -5538.1362
